# Scheduled runner update: refresh market-board price snapshots and
# recompute the dependent Leve profit columns across all Job sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2936

$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2468

$ws.Range("H113").Value = 1480
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -8008

$ws.Range("H138").Value = 4311941
$ws.Range("I138").Value = 1153.5862
$ws.Range("J138").Value = 8622729
$ws.Range("K138").Value = 3460.7586
$ws.Range("L138").Value = 25868187
$ws.Range("M138").Value = 1679.2414
$ws.Range("N138").Value = -25878467

# ---------------------------------------------------------------
# ARM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15169.517
$ws.Range("I32").Value = 19077.227
$ws.Range("J32").Value = 5617.3335
$ws.Range("K32").Value = 19077.227
$ws.Range("L32").Value = 5617.3335
$ws.Range("M32").Value = -18790.227
$ws.Range("N32").Value = -6191.3335

$ws.Range("H42").Value = 1028
$ws.Range("I42").Value = 1028
$ws.Range("K42").Value = 1028
$ws.Range("M42").Value = -542

# ---------------------------------------------------------------
# BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9000
$ws.Range("I19").Value = 9000
$ws.Range("K19").Value = 9000
$ws.Range("M19").Value = -8827

$ws.Range("H134").Value = 4466.857
$ws.Range("I134").Value = 2957.818
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 8873.454000000002
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -6338.454000000002
$ws.Range("N134").Value = -35070

# ---------------------------------------------------------------
# CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2970.6667
$ws.Range("I11").Value = 1900
$ws.Range("J11").Value = 3506
$ws.Range("K11").Value = 1900
$ws.Range("L11").Value = 3506
$ws.Range("M11").Value = -1760
$ws.Range("N11").Value = -3786

$ws.Range("H13").Value = 3005
$ws.Range("J13").Value = 3005
$ws.Range("L13").Value = 3005
$ws.Range("N13").Value = -3283

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H95").Value = 14724.6
$ws.Range("J95").Value = 14724.6
$ws.Range("L95").Value = 14724.6
$ws.Range("N95").Value = -20216.6

$ws.Range("H97").Value = 21087.3
$ws.Range("I97").Value = 20100
$ws.Range("J97").Value = 21197
$ws.Range("K97").Value = 20100
$ws.Range("L97").Value = 21197
$ws.Range("M97").Value = -19109
$ws.Range("N97").Value = -23179

$ws.Range("H102").Value = 24000
$ws.Range("J102").Value = 24000
$ws.Range("L102").Value = 24000
$ws.Range("N102").Value = -28868

# ---------------------------------------------------------------
# CUL
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H134").Value = 3981.1924
$ws.Range("I134").Value = 2581.9375
$ws.Range("J134").Value = 6220
$ws.Range("K134").Value = 7745.8125
$ws.Range("L134").Value = 18660
$ws.Range("M134").Value = -2675.8125
$ws.Range("N134").Value = -28800

$ws.Range("H137").Value = 4371.7036
$ws.Range("I137").Value = 3467.0908
$ws.Range("J137").Value = 4993.625
$ws.Range("K137").Value = 10401.2724
$ws.Range("L137").Value = 14980.875
$ws.Range("M137").Value = -5301.2724
$ws.Range("N137").Value = -25180.875

# ---------------------------------------------------------------
# GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2550
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -5224

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H102").Value = 3100.95
$ws.Range("I102").Value = 3432.353
$ws.Range("J102").Value = 1223
$ws.Range("K102").Value = 3432.353
$ws.Range("L102").Value = 1223
$ws.Range("M102").Value = -1810.353
$ws.Range("N102").Value = -4467

# ---------------------------------------------------------------
# LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5152.304
$ws.Range("I7").Value = 5042.7144
$ws.Range("K7").Value = 5042.7144
$ws.Range("M7").Value = -4930.7144

$ws.Range("H35").Value = 4134
$ws.Range("I35").Value = 1900
$ws.Range("J35").Value = 4692.5
$ws.Range("K35").Value = 1900
$ws.Range("L35").Value = 4692.5
$ws.Range("M35").Value = -1564
$ws.Range("N35").Value = -5364.5

$ws.Range("H39").Value = 240000
$ws.Range("I39").Value = 240000
$ws.Range("K39").Value = 240000
$ws.Range("M39").Value = -239540

$ws.Range("H40").Value = 6963.143
$ws.Range("I40").Value = 6728.4
$ws.Range("J40").Value = 7550
$ws.Range("K40").Value = 6728.4
$ws.Range("L40").Value = 7550
$ws.Range("M40").Value = -6592.4
$ws.Range("N40").Value = -7822

$ws.Range("H56").Value = 25600
$ws.Range("J56").Value = 31400
$ws.Range("L56").Value = 31400
$ws.Range("N56").Value = -32782

$ws.Range("H92").Value = 12000
$ws.Range("J92").Value = 12000
$ws.Range("L92").Value = 12000
$ws.Range("N92").Value = -16992

$ws.Range("H94").Value = 52165
$ws.Range("J94").Value = 52165
$ws.Range("L94").Value = 52165
$ws.Range("N94").Value = -53517

$ws.Range("H126").Value = 5152.304
$ws.Range("I126").Value = 5042.7144
$ws.Range("K126").Value = 15128.1432
$ws.Range("M126").Value = -12658.1432

# ---------------------------------------------------------------
# WVR
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 25000
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25406

$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -21108

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H126").Value = 3466.75
$ws.Range("I126").Value = 1963.7273
$ws.Range("K126").Value = 5891.1819
$ws.Range("M126").Value = -3421.1819
